$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H43").Value = 1059.75
$ws.Range("J43").Value = 1989
$ws.Range("L43").Value = 1989
$ws.Range("N43").Value = -2127
$ws.Range("H53").Value = 360.18182
$ws.Range("I53").Value = 345.125
$ws.Range("K53").Value = 345.125
$ws.Range("M53").Value = 291.875
$ws.Range("H64").Value = 4985.4287
$ws.Range("J64").Value = 4900
$ws.Range("L64").Value = 4900
$ws.Range("N64").Value = -5396
$ws.Range("H67").Value = 4985.4287
$ws.Range("J67").Value = 4900
$ws.Range("L67").Value = 4900
$ws.Range("N67").Value = -6616
$ws.Range("H70").Value = 2630.4614
$ws.Range("I70").Value = 1032.8334
$ws.Range("J70").Value = 3999.8572
$ws.Range("K70").Value = 3098.5002
$ws.Range("L70").Value = 11999.5716
$ws.Range("M70").Value = -2828.5002
$ws.Range("N70").Value = -12539.5716
$ws.Range("H73").Value = 2630.4614
$ws.Range("I73").Value = 1032.8334
$ws.Range("J73").Value = 3999.8572
$ws.Range("K73").Value = 3098.5002
$ws.Range("L73").Value = 11999.5716
$ws.Range("M73").Value = -2162.5002
$ws.Range("N73").Value = -13871.5716
$ws.Range("H74").Value = 86667
$ws.Range("I74").Value = 4000.4
$ws.Range("K74").Value = 4000.4
$ws.Range("M74").Value = -3064.4
$ws.Range("H77").Value = 86667
$ws.Range("I77").Value = 4000.4
$ws.Range("K77").Value = 20002
$ws.Range("M77").Value = -15322
$ws.Range("H98").Value = 909.2632
$ws.Range("I98").Value = 499
$ws.Range("J98").Value = 1798.1666
$ws.Range("K98").Value = 499
$ws.Range("L98").Value = 1798.1666
$ws.Range("M98").Value = 999
$ws.Range("N98").Value = -4794.1666
$ws.Range("H112").Value = 1324.75
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H113").Value = 3369.7144
$ws.Range("I113").Value = 3477.6
$ws.Range("J113").Value = 3100
$ws.Range("K113").Value = 3477.6
$ws.Range("L113").Value = 3100
$ws.Range("M113").Value = -223.5999999999999
$ws.Range("N113").Value = -9608
$ws.Range("H122").Value = 909.2632
$ws.Range("I122").Value = 499
$ws.Range("J122").Value = 1798.1666
$ws.Range("K122").Value = 1497
$ws.Range("L122").Value = 5394.4998
$ws.Range("M122").Value = 953
$ws.Range("N122").Value = -10294.4998
$ws.Range("H125").Value = 4674.4
$ws.Range("I125").Value = 2855.625
$ws.Range("J125").Value = 11949.5
$ws.Range("K125").Value = 25700.625
$ws.Range("L125").Value = 107545.5
$ws.Range("M125").Value = -23240.625
$ws.Range("N125").Value = -112465.5
$ws.Range("H127").Value = 1011.61536
$ws.Range("I127").Value = 512.5833
$ws.Range("K127").Value = 1537.7499
$ws.Range("M127").Value = 3422.2501
$ws.Range("H138").Value = 3363.6365
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("H141").Value = 4365.8
$ws.Range("I141").Value = 4365.8
$ws.Range("K141").Value = 13097.4
$ws.Range("M141").Value = -7917.400000000001

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H61").Value = 2537
$ws.Range("J61").Value = 3494.5
$ws.Range("L61").Value = 3494.5
$ws.Range("N61").Value = -3918.5
$ws.Range("H74").Value = 5278.8
$ws.Range("I74").Value = 4752.6665
$ws.Range("J74").Value = 10014
$ws.Range("K74").Value = 4752.6665
$ws.Range("L74").Value = 10014
$ws.Range("M74").Value = -3878.6665
$ws.Range("N74").Value = -11762
$ws.Range("H77").Value = 5278.8
$ws.Range("I77").Value = 4752.6665
$ws.Range("J77").Value = 10014
$ws.Range("K77").Value = 23763.3325
$ws.Range("L77").Value = 50070
$ws.Range("M77").Value = -19395.3325
$ws.Range("N77").Value = -58806
$ws.Range("H132").Value = 1902.6364
$ws.Range("I132").Value = 1898
$ws.Range("K132").Value = 5694
$ws.Range("M132").Value = -3164
$ws.Range("H136").Value = 2537
$ws.Range("J136").Value = 3494.5
$ws.Range("L136").Value = 10483.5
$ws.Range("N136").Value = -15583.5

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H134").Value = 3231.1904
$ws.Range("I134").Value = 3231.1904
$ws.Range("K134").Value = 9693.5712
$ws.Range("M134").Value = -7158.5712

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H16").Value = 1197.8
$ws.Range("I16").Value = 1197.8
$ws.Range("K16").Value = 1197.8
$ws.Range("M16").Value = -910.8
$ws.Range("H58").Value = 1971.25
$ws.Range("I58").Value = 1971.25
$ws.Range("K58").Value = 1971.25
$ws.Range("M58").Value = -1768.25
$ws.Range("H94").Value = 1891.625
$ws.Range("I94").Value = 1955.6666
$ws.Range("J94").Value = 1699.5
$ws.Range("K94").Value = 1955.6666
$ws.Range("L94").Value = 1699.5
$ws.Range("M94").Value = -1504.6666
$ws.Range("N94").Value = -2601.5
$ws.Range("H99").Value = 6006.125
$ws.Range("I99").Value = 5924.8335
$ws.Range("K99").Value = 5924.8335
$ws.Range("M99").Value = -4426.8335
$ws.Range("H112").Value = 69999.5
$ws.Range("J112").Value = 69999.5
$ws.Range("L112").Value = 69999.5
$ws.Range("N112").Value = -72953.5
$ws.Range("H113").Value = 1197.8
$ws.Range("I113").Value = 1197.8
$ws.Range("K113").Value = 1197.8
$ws.Range("M113").Value = 972.2
$ws.Range("H126").Value = 6006.125
$ws.Range("I126").Value = 5924.8335
$ws.Range("K126").Value = 17774.5005
$ws.Range("M126").Value = -15304.5005
$ws.Range("H134").Value = 1762
$ws.Range("J134").Value = 949.5
$ws.Range("L134").Value = 2848.5
$ws.Range("N134").Value = -7918.5
$ws.Range("H136").Value = 1971.25
$ws.Range("I136").Value = 1971.25
$ws.Range("K136").Value = 5913.75
$ws.Range("M136").Value = -3363.75

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H5").Value = 2869.5715
$ws.Range("J5").Value = 2122
$ws.Range("L5").Value = 6366
$ws.Range("N5").Value = -6590
$ws.Range("H14").Value = 37647.375
$ws.Range("I14").Value = 37647.375
$ws.Range("K14").Value = 112942.125
$ws.Range("M14").Value = -112769.125
$ws.Range("H117").Value = 755.6667
$ws.Range("I117").Value = 759.5
$ws.Range("J117").Value = 748
$ws.Range("K117").Value = 2278.5
$ws.Range("L117").Value = 2244
$ws.Range("M117").Value = 1163.5
$ws.Range("N117").Value = -9128
$ws.Range("H128").Value = 339792.88
$ws.Range("I128").Value = 339792.88
$ws.Range("K128").Value = 1019378.64
$ws.Range("M128").Value = -1014398.64
$ws.Range("H135").Value = 2869.5715
$ws.Range("J135").Value = 2122
$ws.Range("L135").Value = 19098
$ws.Range("N135").Value = -24168

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H43").Value = 5550
$ws.Range("H57").Value = 22600
$ws.Range("I57").Value = 12000
$ws.Range("K57").Value = 12000
$ws.Range("M57").Value = -11180
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("H107").Value = 975
$ws.Range("I107").Value = 300
$ws.Range("K107").Value = 300
$ws.Range("M107").Value = 1620
$ws.Range("H113").Value = 3504.4
$ws.Range("I113").Value = 3504.4
$ws.Range("K113").Value = 3504.4
$ws.Range("M113").Value = -1334.4
$ws.Range("H126").Value = 2847.2144
$ws.Range("I126").Value = 2488.4167
$ws.Range("K126").Value = 7465.250100000001
$ws.Range("M126").Value = -4995.250100000001
$ws.Range("H132").Value = 4499.6665
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 6000
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H40").Value = 4133.3335
$ws.Range("I40").Value = 5200
$ws.Range("K40").Value = 5200
$ws.Range("M40").Value = -5064
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 2693.7693
$ws.Range("I122").Value = 2613.2222
$ws.Range("J122").Value = 2875
$ws.Range("K122").Value = 7839.6666
$ws.Range("L122").Value = 8625
$ws.Range("M122").Value = -5389.6666
$ws.Range("N122").Value = -13525
$ws.Range("H132").Value = 6043.4165
$ws.Range("I132").Value = 2701
$ws.Range("J132").Value = 7714.625
$ws.Range("K132").Value = 8103
$ws.Range("L132").Value = 23143.875
$ws.Range("M132").Value = -5573
$ws.Range("N132").Value = -28203.875

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H107").Value = 398
$ws.Range("I107").Value = 398
$ws.Range("K107").Value = 1194
$ws.Range("M107").Value = 726
$ws.Range("H132").Value = 1881.7142
$ws.Range("I132").Value = 1881.7142
$ws.Range("K132").Value = 5645.142599999999
$ws.Range("M132").Value = -3115.142599999999
$ws.Range("H136").Value = 6486.5
$ws.Range("I136").Value = 5315.3335
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 15946.0005
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -13396.0005
$ws.Range("N136").Value = -35100
